# Fix AMI IDs for Ubuntu 14 on AWS (Sheet1, row 4)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 4 is "Ubuntu 14":
#   D4 (us-east-1)    : ami-2d57433a -> ami-cfa100d9
#   F4 (us-west-1)    : ami-e7277687 -> ami-29752c49
#   G4 (us-west-2)    : ami-900ebaf0 -> removed (cleared)
#   H4 (eu-west-1)    : ami-75c9e906 -> ami-32517b54
#   I4 (eu-central-1) : ami-162ded79 -> ami-0738ec68

$ws.Range("D4").Value = "ami-cfa100d9"
$ws.Range("F4").Value = "ami-29752c49"
$ws.Range("G4").ClearContents()
$ws.Range("H4").Value = "ami-32517b54"
$ws.Range("I4").Value = "ami-0738ec68"

# Update the active selection cursor to J10 as recorded in the saved file
$ws.Range("J10").Select()
